$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-10-25 02:20:06"

$wsZhCn.Range("H3").Value = "2016-10-25 02:19:54"
$wsZhCn.Range("K3").Value = "2016-10-25 02:20:37"

$wsDeDe.Range("K3").Value = "2016-10-25 02:20:55"
